# The commit swaps the contents of ppt/theme/theme1.xml and ppt/theme/theme2.xml:
#   before: theme1.xml = "Office Theme" palette, theme2.xml = "Integral" palette
#   after : theme1.xml = "Integral" palette,     theme2.xml = "Office Theme" palette
#
# theme2.xml is the presentation's live/active theme (the one bound to the
# slide master and therefore to every slide), and it is the only theme part
# reachable through the PowerPoint object model. We reassign its twelve
# scheme colors, in place, to the "Office Theme" palette so the deck's
# visible design now matches what the diff shows theme2.xml containing
# after the edit.
#
# ThemeColorScheme.Colors(i).RGB uses the standard COM COLORREF (BGR) byte
# order, so each target 0xRRGGBB value below is supplied as 0xBBGGRR.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$scheme = $theme.ThemeColorScheme

$scheme.Colors(1).RGB  = 0x000000   # dk1
$scheme.Colors(2).RGB  = 0xFFFFFF   # lt1
$scheme.Colors(3).RGB  = 0x6A5444   # dk2
$scheme.Colors(4).RGB  = 0xE6E6E7   # lt2
$scheme.Colors(5).RGB  = 0xD59B5B   # accent1
$scheme.Colors(6).RGB  = 0x317DED   # accent2
$scheme.Colors(7).RGB  = 0xA5A5A5   # accent3
$scheme.Colors(8).RGB  = 0x00C0FF   # accent4
$scheme.Colors(9).RGB  = 0xC47244   # accent5
$scheme.Colors(10).RGB = 0x47AD70   # accent6
$scheme.Colors(11).RGB = 0xC16305   # hlink
$scheme.Colors(12).RGB = 0x724F95   # folHlink
